$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The leftmost column (A) held a stray duplicate of the GENE number and was
# styled like a header cell by mistake. Drop it entirely so every later
# column shifts one slot to the left: old B->A (QS_Astral25), C->B
# (FNRATE_ASTRAL), D->C (TAXON), E->D (MODEL_CONDITION), F->E (GENE).
$ws.Range("A1").EntireColumn.Delete()
